$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-9, columns A:T
# (Sending cluster, Ligand symbol, Receptor symbol, Target cluster, then 16 numeric metrics)
$data = @(
    @("ECs", "Cd48", "Cd2", "M1", 1, 0.3333333333333333, 0.06575866666666667, 0.197276,            0.0003210775944521127, 0.000321193332710493, 3, 1, 0.7139163333333333, 2.141749, 0.5325697612629201, 0.5325697612629202, 0.04694618619155556, 0.422515675724,      0.0001709962178242343, 0.0001710578565208689),
    @("ECs", "Cd48", "Cd2", "M2", 1, 0.3333333333333333, 0.06575866666666667, 0.197276,            0.0003210775944521127, 0.000321193332710493, 3, 1, 0.6265959999999999, 1.879788, 0.4674302387370798, 0.4674302387370799, 0.04120411749866667, 0.370837057488,      0.0001500813766278783, 0.0001501354761896241),
    @("M1",  "Cd48", "Cd2", "M1", 3, 1,                  93.51811733333334,  280.554352,           0.4566177155519742,    0.4567823117120824,    3, 1, 0.7139163333333333, 2.141749, 0.5325697612629201, 0.5325697612629202, 66.76411142684978,   600.8770028416479,  0.2431807877599348,    0.2432684466976285),
    @("M1",  "Cd48", "Cd2", "M2", 3, 1,                  93.51811733333334,  280.554352,           0.4566177155519742,    0.4567823117120824,    3, 1, 0.6265959999999999, 1.879788, 0.4674302387370798, 0.4674302387370799, 58.59807824859733,   527.3827042373759,  0.2134369277920393,    0.2135138650144539),
    @("M2",  "Cd48", "Cd2", "M1", 3, 1,                  111.0008776666667,  333.002633,           0.5419801919638461,    0.5421755585810701,    3, 1, 0.7139163333333333, 2.141749, 0.5325697612629201, 0.5325697612629202, 79.24533958056855,   713.2080562251169,  0.2886422614434171,    0.2887463077961109),
    @("M2",  "Cd48", "Cd2", "M2", 3, 1,                  111.0008776666667,  333.002633,           0.5419801919638461,    0.5421755585810701,    3, 1, 0.6265959999999999, 1.879788, 0.4674302387370798, 0.4674302387370799, 69.55270594242266,   625.974353481804,   0.2533379305204289,    0.2534292507849593),
    @("sCs", "Cd48", "Cd2", "M1", 1, 0.5,                0.2213985,          0.442797,             0.001081014889727682,  0.0007209363741367839, 3, 1, 0.7139163333333333, 2.141749, 0.5325697612629201, 0.5325697612629202, 0.1580600053255,     0.9483600319529999, 0.0005757158417439334, 0.0003839489126597823),
    @("sCs", "Cd48", "Cd2", "M2", 1, 0.5,                0.2213985,          0.442797,             0.001081014889727682,  0.0007209363741367839, 3, 1, 0.6265959999999999, 1.879788, 0.4674302387370798, 0.4674302387370799, 0.138727414506,      0.8323644870359999, 0.0005052990479837483, 0.0003369874614770016)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
